# Natmi following Dr Hou advice
# Update the LR-pair table (columns E..T) for rows 2-10 with recomputed
# ligand/receptor expression statistics after the ligand/receptor
# expressing-cell threshold changed from 1 to 3 expressing cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hashtable of new values per row; keys are column letters, values are the
# recalculated statistics taken from the updated natmi run.
$data = @{
    2  = @{ E=3; G=12.07002266666666;  H=36.21006799999999; I=0.7601982364861632; J=0.7601982364861634; K=3; M=7.214110666666667;  N=21.642332; O=0.4688823795981188; P=0.4688823795981188; Q=87.07447926650843;  R=783.6703133985758;  S=0.3564435580899257;  T=0.3564435580899258  }
    3  = @{ E=3; G=12.07002266666666;  H=36.21006799999999; I=0.7601982364861632; J=0.7601982364861634; K=3; M=7.110350666666666;  N=21.331052; O=0.4621384803214003; P=0.4621384803214003; Q=85.82209371461509;  R=772.3988434315359;  S=0.351316857752724;   T=0.351316857752724   }
    4  = @{ E=3; G=12.07002266666666;  H=36.21006799999999; I=0.7601982364861632; J=0.7601982364861634; K=3; M=1.061296333333333;  N=3.183889;  O=0.06897914008048092; P=0.06897914008048092; Q=12.80987079938355;  R=115.288837194452;   S=0.05243782064351362; T=0.05243782064351363 }
    5  = @{ E=3; G=1.308268;           H=3.924804;           I=0.08239777620284613;J=0.08239777620284613;K=3; M=7.214110666666667;  N=21.642332; O=0.4688823795981188; P=0.4688823795981188; Q=9.437990133658666;  R=84.941911202928;    S=0.03863486537958374; T=0.03863486537958374 }
    6  = @{ E=3; G=1.308268;           H=3.924804;           I=0.08239777620284613;J=0.08239777620284613;K=3; M=7.110350666666666;  N=21.331052; O=0.4621384803214003; P=0.4621384803214003; Q=9.302244245978667;  R=83.720198213808;    S=0.03807918307624616; T=0.03807918307624616 }
    7  = @{ E=3; G=1.308268;           H=3.924804;           I=0.08239777620284613;J=0.08239777620284613;K=3; M=1.061296333333333;  N=3.183889;  O=0.06897914008048092; P=0.06897914008048092; Q=1.388460031417333;  R=12.496140282756;    S=0.005683727747016241;T=0.005683727747016241}
    8  = @{ E=3; G=2.499176666666667;  H=7.49753;            I=0.1574039873109905; J=0.1574039873109906; K=3; M=7.214110666666667;  N=21.642332; O=0.4688823795981188; P=0.4688823795981188; Q=18.02933704888444;  R=162.26403343996;    S=0.07380395612860935; T=0.07380395612860936 }
    9  = @{ E=3; G=2.499176666666667;  H=7.49753;            I=0.1574039873109905; J=0.1574039873109906; K=3; M=7.110350666666666;  N=21.331052; O=0.4621384803214003; P=0.4621384803214003; Q=17.77002247795111;  R=159.93020230156;    S=0.07274243949243014; T=0.07274243949243016 }
    10 = @{ E=3; G=2.499176666666667;  H=7.49753;            I=0.1574039873109905; J=0.1574039873109906; K=3; M=1.061296333333333;  N=3.183889;  O=0.06897914008048092; P=0.06897914008048092; Q=2.652367032685555;  R=23.87130329417;     S=0.01085759168995106; T=0.01085759168995106 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
